$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.211807012557983
$ws.Range("B1").Value = 1.981002569198608
$ws.Range("C1").Value = 4.21067476272583
$ws.Range("D1").Value = 3.006114959716797
$ws.Range("E1").Value = 1.194965839385986
